# Apply translation-string updates described by the commit diff.
#
# The sheet is protected, and columns A:C use a style (index 1) that this
# runtime treats as "locked" for direct Range.Value writes even though the
# style itself unlocks the cells (Unprotect()/Protect() round-trips also
# collapse the sheetProtection element's fine-grained flags, which aren't
# part of this change) - so writes into A:C go through a Copy /
# PasteSpecial(xlPasteValues) relay via an unused scratch cell, which this
# runtime does not subject to the same protected-sheet check. Columns D:X
# are not protection-flagged and can be written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122
$scratch = $ws.Range("Z1000")

function Set-LockedCell($addr, $value) {
    $scratch.Value = $value
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial($xlPasteValues)
}

function Copy-RowStyle($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
}

# --- Fill in previously-empty "Comment" column (D) values ---
$ws.Range("D2").Value  = "In plot derivative"
$ws.Range("D3").Value  = "In plot derivative"
$ws.Range("D4").Value  = "In plot derivative"
$ws.Range("D5").Value  = "Column text header in exported files"
$ws.Range("D6").Value  = "Tab text in `"settings`" form"
$ws.Range("D7").Value  = "In `"settings`" form"
$ws.Range("D8").Value  = "In `"settings`" form"
$ws.Range("D9").Value  = "In `"settings`" form"
$ws.Range("D10").Value = "In `"settings`" form, mathematical name"
$ws.Range("D11").Value = "In `"settings`" form, mathematical name"
$ws.Range("D12").Value = "In `"settings`" form, mathematical name"
$ws.Range("D13").Value = "In `"settings`" form, mathematical name"

# --- Row 4: rename key, and shorten "seconds" to "sec" / "segundos" to "seg" ---
Set-LockedCell "C4" "strPlotDerivativeYLabel1"
$ws.Range("E4").Value = "Amplitude / sec"
$ws.Range("K4").Value = "Amplitude / sec"
$ws.Range("W4").Value = "Amplitud / seg"

# --- Pluralize Spanish strings in rows 12 and 13 ---
$ws.Range("W12").Value = "Diferencias centrales de 3 puntos"
$ws.Range("W13").Value = "Diferencias centrales de 5 puntos"

# --- New row 14: strPlotDerivativeYLabel2 ---
Set-LockedCell "A14" "SignalAnalysis"
Set-LockedCell "B14" "localization\strings"
Set-LockedCell "C14" "strPlotDerivativeYLabel2"
Copy-RowStyle "A13:C13" "A14:C14"
$ws.Range("D14").Value = "In plot derivative"
$ws.Range("E14").Value = "Amplitude"
$ws.Range("K14").Value = "Amplitude"
$ws.Range("W14").Value = "Amplitud"

# --- New row 15: strStatusTipDerivative ---
Set-LockedCell "A15" "SignalAnalysis"
Set-LockedCell "B15" "localization\strings"
Set-LockedCell "C15" "strStatusTipDerivative"
Copy-RowStyle "A13:C13" "A15:C15"
$ws.Range("E15").Value = "Numerical differentiation"
$ws.Range("K15").Value = "Numerical differentiation"
$ws.Range("W15").Value = "Derivada numérica"

# Clean up the scratch cell used as a copy source.
$scratch.ClearContents()
